$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap Name/Ticker columns (B<->C) and reorder rows, then add new
# header columns D:K with additional metric headers on row 2.

$ws.Range("B2").Value2 = "Name"
$ws.Range("C2").Value2 = "Ticker"
$ws.Range("D2").Value2 = "Price"
$ws.Range("E2").Value2 = "MC"
$ws.Range("F2").Value2 = "EV"
$ws.Range("G2").Value2 = "2022 EV/E"
$ws.Range("H2").Value2 = "2023 EV/E"
$ws.Range("I2").Value2 = "2022 E"
$ws.Range("J2").Value2 = "2023 E"
$ws.Range("K2").Value2 = "2022 RG"

$ws.Range("B3").Value2 = "Berkeley"
$ws.Range("C3").Value2 = "BKG"

$ws.Range("B4").Value2 = "Reckitt Benckiser"
$ws.Range("C4").Value2 = "RKT"

$ws.Range("B5").Value2 = "Aegon NV"
$ws.Range("C5").Value2 = "AEG"

$ws.Range("B6").Value2 = "Forum Energy Technologies"
$ws.Range("C6").Value2 = "FET "

$ws.Range("B7").Value2 = "BP"
$ws.Range("C7").Value2 = "BP."

$ws.Range("B8").Value2 = "Chart Industries"
$ws.Range("C8").Value2 = "GTLS"

# Column B now holds the wider "Name" text and needs the autofit width
# that used to live on column C; column C (tickers) reverts to default.
$ws.Columns("C").ColumnWidth = 8.43
$ws.Columns("B").ColumnWidth = 23.140625

$ws.Range("J25").Select()
